$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row above row 833, pushing the existing row 833
# (and everything below it) down by one. This matches the diff: a new
# data point for 2026/02/22 gets inserted into the middle of the table,
# and every subsequent row shifts down one index (dimension A1:D874 -> A1:D875).
$ws.Rows.Item(833).EntireRow.Insert()

# Fill in the newly inserted row. Column A holds a date string formatted
# as "yyyy/mm/dd" (e.g. "2026/12/29") but stored as literal text in the
# original workbook (t="inlineStr"), not as a real date serial. Force the
# cell to text first so Excel's auto date-detection doesn't convert
# "2026/02/22" into a date value, then restore the default "Normal" style
# so the cell doesn't end up carrying a stray explicit number format.
$ws.Cells.Item(833, 1).NumberFormat = "@"
$ws.Cells.Item(833, 1).Value = "2026/02/22"
$ws.Cells.Item(833, 1).Style = "Normal"

$ws.Cells.Item(833, 2).Value = "日"
$ws.Cells.Item(833, 3).Value = 13
$ws.Cells.Item(833, 4).Value = 47

Write-Output "Inserted row 833 (2026/02/22, 日, 13, 47); table now A1:D875"
